# Insert a new weekly price record for Albahaca at Terminal La Palmera de La Serena.
# This inserts a new row at row 36, shifting all existing data rows (36-143) down by one
# (becoming rows 37-144), and populates the new row 36 with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 36 - existing row 36 (and below) shift down to 37.
$ws.Rows("36").Insert()

# Fill in the constant columns (same value across the whole dataset).
$ws.Range("A36").Value = 8
$ws.Range("B36").Value = "Terminal La Palmera de La Serena"
$ws.Range("C36").Value = "Coquimbo"
$ws.Range("D36").Value = 44910
$ws.Range("E36").Value = 4
$ws.Range("F36").Value = 100112052
$ws.Range("G36").Value = "Albahaca"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 1000
$ws.Range("K36").Value = 3800
$ws.Range("L36").Value = 4000
$ws.Range("M36").Value = 3900
$ws.Range("N36").Value = "$/paquete"
$ws.Range("O36").Value = "Región de Arica y Parinacota"
$ws.Range("P36").Value = 3900
$ws.Range("Q36").Value = 1
$ws.Range("R36").Value = "Hortaliza"

# Match the date-cell style/number format used by the rest of column D.
$ws.Range("D36").NumberFormat = $ws.Range("D37").NumberFormat
